# This edit re-shuffles the content of rows 5-25 (except rows 8, 10 and 20,
# which stay as-is) in the "Artfynd" sheet: each row keeps its row number /
# formatting, but the field values that used to live in one row now live in
# another row, per the mapping below (target row -> source row, both
# referring to the *original* row numbering before any changes were made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (content originally found in the source row ends
# up in the target row)
$mapping = [ordered]@{
    5  = 6
    6  = 7
    7  = 5
    9  = 23
    11 = 17
    12 = 9
    13 = 19
    14 = 12
    15 = 21
    16 = 25
    17 = 13
    18 = 15
    19 = 16
    21 = 11
    22 = 24
    23 = 22
    24 = 18
    25 = 14
}

# First capture the full row contents (columns A:AY) for every row that is
# either a source or a target of the permutation, before any writes happen.
$snapshot = @{}
foreach ($srcRow in $mapping.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $snapshot[$srcRow] = $ws.Range("A" + $srcRow + ":AY" + $srcRow).Value()
    }
}

# Columns I (Antal), Y (Startdatum) and AA (Slutdatum) hold values that look
# like numbers / dates but are actually stored as plain text in the source
# data. Excel's automatic type detection would otherwise silently convert
# them to numbers / real dates when the values are written back through the
# object model, so force those columns to Text format first to keep them as
# plain strings.
foreach ($targetRow in $mapping.Keys) {
    $ws.Range("I" + $targetRow).NumberFormat = "@"
    $ws.Range("Y" + $targetRow).NumberFormat = "@"
    $ws.Range("AA" + $targetRow).NumberFormat = "@"
}

# Now write each captured row into its target row.
foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $ws.Range("A" + $targetRow + ":AY" + $targetRow).Value = $snapshot[$srcRow]
}
